# Generate Report for Archive
# The localization status moved from "Ready for handoff" to "In Translation"
# for every tracked file, on every sheet (Overview + one sheet per locale).
# Shortening that status text lets Excel's column autosizing shrink the two
# "status" columns on each sheet, so we re-apply the (smaller) column widths
# after updating the text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# --- Per-locale detail sheets: Status column (C) for both rows ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Re-fit the status columns now that the text is shorter ---
# (target display width ~= 13.41 characters; 12.5 is the COM ColumnWidth
# input that lands on the nearest width this host can represent)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
